# Insert a new data row at row 182 (shifting existing rows 182.. down by one)
# and populate it with the new record described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 182; this shifts rows 182..282 down to 183..282
$ws.Rows.Item(182).Insert()

# Populate the newly inserted row 182 with the new record's values.
$ws.Cells.Item(182, 1).Value = 10
$ws.Cells.Item(182, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(182, 3).Value = "La Araucanía"
$ws.Cells.Item(182, 4).Value = 44438
$ws.Cells.Item(182, 5).Value = 9
$ws.Cells.Item(182, 6).Value = "Fruta"
$ws.Cells.Item(182, 7).Value = 100101
$ws.Cells.Item(182, 8).Value = "Berries"
$ws.Cells.Item(182, 9).Value = 100101007
$ws.Cells.Item(182, 10).Value = "Kiwi"
$ws.Cells.Item(182, 11).Value = "Hayward"
$ws.Cells.Item(182, 12).Value = "Primera"
$ws.Cells.Item(182, 13).Value = 235
$ws.Cells.Item(182, 14).Value = 18000
$ws.Cells.Item(182, 15).Value = 19000
$ws.Cells.Item(182, 16).Value = 18468
$ws.Cells.Item(182, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(182, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(182, 19).Value = 1026
$ws.Cells.Item(182, 20).Value = 18
